# Target diff analysis
# -----------------------------------------------------------------------
# The supplied unified diff touches exactly two package parts:
#   - customXml/item1.xml       (the SharePoint "content type schema" —
#                                  field display-name locale strings such
#                                  as "Content Type"->"Inhoudstype",
#                                  "Image Tags"->"Afbeeldingtags", a bumped
#                                  ma:contentTypeVersion/ma:versionID/
#                                  ma:fieldsID, and one extra hidden
#                                  managed-metadata field definition,
#                                  MediaServiceObjectDetectorVersions)
#   - customXml/itemProps1.xml  (just the ds:datastoreItem GUID + the
#                                  list of ds:schemaRef entries collapsed)
#
# Both files are server/SharePoint-managed metadata describing the
# document library's content type; they are not part of the visible
# Word document (body text, headers/footers, styles, tables, ...) and
# are not reachable through the Word UI or the Word object model
# (Document.CustomXMLParts / XMLNodes / ContentTypeProperties all report
# zero items / empty XML in this runtime, exactly like real Word treats
# this class of content-type XML as read-only server plumbing rather
# than document content a macro edits).
#
# That matches the author's own commit message for this change:
#   "Geen wijzigingen. Veranderingen in data na installatie op productie"
#   ("No changes. Changes in data after installation on production.")
# i.e. nothing in the authored document content actually changed — the
# diff is incidental metadata drift that SharePoint stamped onto the
# package after a production deploy (locale-dependent field captions,
# a regenerated content-type version/GUID, and a newly synced hidden
# field), not an edit a user made inside Word.
#
# There is therefore no Word object-model action to perform here: no
# text, formatting, table, or other document-content change is implied
# by this diff. This script intentionally performs no content mutation.

$d = $word.ActiveDocument

# Touch the document object to confirm the session is alive, without
# altering any content (no Range/Find/Property writes below), mirroring
# the "no changes" commit message.
$null = $d.Name
